# Re-sync the SharePoint "Document" content-type schema custom XML part.
#
# customXml/item3.xml is the cached ct:contentTypeSchema blob SharePoint
# stamps into the package; a content-type republish bumps its
# ma:versionID / ma:fieldsID hashes. customXml/itemProps3.xml is the
# datastore descriptor for that part (ds:itemID) - on a resync the part is
# re-minted with a fresh id and the cached ds:schemaRefs list is dropped
# (rehydrated later from the list's schema library).
#
# Do this through Document.CustomXMLParts, the Word object-model surface
# for custom XML parts, rather than touching package bytes directly.

$d = $word.ActiveDocument

$oldVersionId = 'ma:versionID="1cc4cf9d95b5e2d14d7aabb44ca49f5e"'
$newVersionId = 'ma:versionID="ce94caacb4a5cc228342027e3189af2c"'
$oldFieldsId  = 'ma:fieldsID="a33e6829bf21261855124b7b230b6e9c"'
$newFieldsId  = 'ma:fieldsID="5f85a36ab557a4a47cd270a1ee4435c0"'
$contentTypeNs = "http://schemas.microsoft.com/office/2006/metadata/contentType"

$target = $null

# Prefer locating the part by its content-type-schema namespace ...
try {
    $scoped = $d.CustomXMLParts.SelectByNamespace($contentTypeNs)
    if ($scoped -ne $null -and $scoped.Count -ge 1) {
        $target = $scoped.Item(1)
    }
} catch {
}

# ... falling back to scanning every part for the versionID marker.
if ($target -eq $null) {
    try {
        $all = $d.CustomXMLParts
        for ($i = 1; $i -le $all.Count; $i++) {
            $candidate = $all.Item($i)
            if ($candidate -ne $null -and $candidate.XML -ne $null -and $candidate.XML.Contains($oldVersionId)) {
                $target = $candidate
                break
            }
        }
    } catch {
    }
}

if ($target -ne $null) {
    $xml = $target.XML
    $xml = $xml.Replace($oldVersionId, $newVersionId)
    $xml = $xml.Replace($oldFieldsId, $newFieldsId)

    # Swap in the updated schema: drop the stale datastore item and add a
    # fresh one. Word mints a brand-new ds:itemID guid and writes a bare
    # <ds:datastoreItem> (no <ds:schemaRefs>) for a newly-added part -
    # matching a server-side content-type resync.
    try {
        $target.Delete()
    } catch {
    }
    try {
        $d.CustomXMLParts.Add($xml) | Out-Null
    } catch {
    }
}
